$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update nombre_aides (column C) and montant_total (column D) values for the 2020-05-20 data refresh.
$ws.Cells.Item(2, 3).Value = 23577
$ws.Cells.Item(2, 4).Value = 34424234
$ws.Cells.Item(3, 3).Value = 59445
$ws.Cells.Item(3, 4).Value = 87997202
$ws.Cells.Item(4, 3).Value = 20154
$ws.Cells.Item(4, 4).Value = 30078301
$ws.Cells.Item(5, 3).Value = 5327
$ws.Cells.Item(5, 4).Value = 7971962
$ws.Cells.Item(6, 3).Value = 1054
$ws.Cells.Item(6, 4).Value = 1579391
$ws.Cells.Item(7, 3).Value = 66
$ws.Cells.Item(7, 4).Value = 99000
$ws.Cells.Item(10, 3).Value = 25229
$ws.Cells.Item(10, 4).Value = 34872548
$ws.Cells.Item(11, 3).Value = 6186
$ws.Cells.Item(11, 4).Value = 9052541
$ws.Cells.Item(12, 3).Value = 17459
$ws.Cells.Item(12, 4).Value = 25829560
$ws.Cells.Item(13, 3).Value = 5425
$ws.Cells.Item(13, 4).Value = 8101348
$ws.Cells.Item(14, 3).Value = 1280
$ws.Cells.Item(14, 4).Value = 1915995
$ws.Cells.Item(15, 3).Value = 232
$ws.Cells.Item(15, 4).Value = 345266
$ws.Cells.Item(17, 3).Value = 6343
$ws.Cells.Item(17, 4).Value = 8597993
$ws.Cells.Item(18, 3).Value = 8636
$ws.Cells.Item(18, 4).Value = 12589224
$ws.Cells.Item(19, 3).Value = 21391
$ws.Cells.Item(19, 4).Value = 31663560
$ws.Cells.Item(20, 3).Value = 6779
$ws.Cells.Item(20, 4).Value = 10133340
$ws.Cells.Item(21, 3).Value = 1629
$ws.Cells.Item(21, 4).Value = 2438302
$ws.Cells.Item(22, 3).Value = 241
$ws.Cells.Item(22, 4).Value = 361124
$ws.Cells.Item(24, 3).Value = 7390
$ws.Cells.Item(24, 4).Value = 10102366
$ws.Cells.Item(25, 3).Value = 4893
$ws.Cells.Item(25, 4).Value = 7139630
$ws.Cells.Item(26, 3).Value = 15061
$ws.Cells.Item(26, 4).Value = 22283690
$ws.Cells.Item(27, 3).Value = 5118
$ws.Cells.Item(27, 4).Value = 7652024
$ws.Cells.Item(28, 3).Value = 1233
$ws.Cells.Item(28, 4).Value = 1848991
$ws.Cells.Item(29, 3).Value = 185
$ws.Cells.Item(29, 4).Value = 277500
$ws.Cells.Item(31, 3).Value = 5212
$ws.Cells.Item(31, 4).Value = 7006788
$ws.Cells.Item(32, 3).Value = 1724
$ws.Cells.Item(32, 4).Value = 2487621
$ws.Cells.Item(33, 3).Value = 4553
$ws.Cells.Item(33, 4).Value = 6690376
$ws.Cells.Item(34, 3).Value = 1841
$ws.Cells.Item(34, 4).Value = 2740922
$ws.Cells.Item(35, 3).Value = 471
$ws.Cells.Item(35, 4).Value = 703041
$ws.Cells.Item(36, 3).Value = 90
$ws.Cells.Item(36, 4).Value = 135000
$ws.Cells.Item(38, 3).Value = 1163
$ws.Cells.Item(38, 4).Value = 1581648
$ws.Cells.Item(39, 3).Value = 10966
$ws.Cells.Item(39, 4).Value = 15998013
$ws.Cells.Item(40, 3).Value = 33839
$ws.Cells.Item(40, 4).Value = 50046730
$ws.Cells.Item(41, 3).Value = 12517
$ws.Cells.Item(41, 4).Value = 18704934
$ws.Cells.Item(42, 3).Value = 3461
$ws.Cells.Item(42, 4).Value = 5181870
$ws.Cells.Item(43, 3).Value = 601
$ws.Cells.Item(43, 4).Value = 900436
$ws.Cells.Item(44, 3).Value = 35
$ws.Cells.Item(44, 4).Value = 51691
$ws.Cells.Item(46, 3).Value = 10331
$ws.Cells.Item(46, 4).Value = 14076349
$ws.Cells.Item(47, 3).Value = 984
$ws.Cells.Item(47, 4).Value = 1424097
$ws.Cells.Item(48, 3).Value = 3685
$ws.Cells.Item(48, 4).Value = 5435059
$ws.Cells.Item(49, 3).Value = 1388
$ws.Cells.Item(49, 4).Value = 2074964
$ws.Cells.Item(50, 3).Value = 423
$ws.Cells.Item(50, 4).Value = 632000
$ws.Cells.Item(51, 3).Value = 89
$ws.Cells.Item(51, 4).Value = 133500
$ws.Cells.Item(52, 3).Value = 2376
$ws.Cells.Item(52, 4).Value = 3306602
$ws.Cells.Item(53, 3).Value = 349
$ws.Cells.Item(53, 4).Value = 507284
$ws.Cells.Item(54, 3).Value = 942
$ws.Cells.Item(54, 4).Value = 1394992
$ws.Cells.Item(55, 3).Value = 381
$ws.Cells.Item(55, 4).Value = 569476
$ws.Cells.Item(56, 3).Value = 129
$ws.Cells.Item(56, 4).Value = 193378
$ws.Cells.Item(58, 3).Value = 451
$ws.Cells.Item(58, 4).Value = 641265
$ws.Cells.Item(59, 3).Value = 9962
$ws.Cells.Item(59, 4).Value = 14476411
$ws.Cells.Item(60, 3).Value = 30343
$ws.Cells.Item(60, 4).Value = 44780815
$ws.Cells.Item(61, 3).Value = 10532
$ws.Cells.Item(61, 4).Value = 15742022
$ws.Cells.Item(62, 3).Value = 2917
$ws.Cells.Item(62, 4).Value = 4366068
$ws.Cells.Item(63, 3).Value = 513
$ws.Cells.Item(63, 4).Value = 769139
$ws.Cells.Item(66, 3).Value = 9853
$ws.Cells.Item(66, 4).Value = 13211945
$ws.Cells.Item(67, 3).Value = 2699
$ws.Cells.Item(67, 4).Value = 3941438
$ws.Cells.Item(68, 3).Value = 7337
$ws.Cells.Item(68, 4).Value = 10828142
$ws.Cells.Item(69, 3).Value = 2600
$ws.Cells.Item(69, 4).Value = 3885113
$ws.Cells.Item(70, 3).Value = 852
$ws.Cells.Item(70, 4).Value = 1276049
$ws.Cells.Item(71, 3).Value = 170
$ws.Cells.Item(71, 4).Value = 253612
$ws.Cells.Item(73, 3).Value = 2828
$ws.Cells.Item(73, 4).Value = 3850533
$ws.Cells.Item(74, 3).Value = 865
$ws.Cells.Item(74, 4).Value = 1254393
$ws.Cells.Item(75, 3).Value = 2979
$ws.Cells.Item(75, 4).Value = 4403846
$ws.Cells.Item(76, 3).Value = 1170
$ws.Cells.Item(76, 4).Value = 1751548
$ws.Cells.Item(77, 3).Value = 405
$ws.Cells.Item(77, 4).Value = 607500
$ws.Cells.Item(78, 3).Value = 84
$ws.Cells.Item(78, 4).Value = 125569
$ws.Cells.Item(80, 3).Value = 1753
$ws.Cells.Item(80, 4).Value = 2358530
$ws.Cells.Item(81, 3).Value = 307
$ws.Cells.Item(81, 4).Value = 456689
$ws.Cells.Item(82, 3).Value = 103
$ws.Cells.Item(82, 4).Value = 154110
$ws.Cells.Item(83, 3).Value = 43
$ws.Cells.Item(83, 4).Value = 64500
$ws.Cells.Item(85, 3).Value = 8
$ws.Cells.Item(85, 4).Value = 12000
$ws.Cells.Item(86, 3).Value = 7007
$ws.Cells.Item(86, 4).Value = 10250101
$ws.Cells.Item(87, 3).Value = 20111
$ws.Cells.Item(87, 4).Value = 29757558
$ws.Cells.Item(88, 3).Value = 6605
$ws.Cells.Item(88, 4).Value = 9872215
$ws.Cells.Item(89, 3).Value = 1748
$ws.Cells.Item(89, 4).Value = 2617655
$ws.Cells.Item(90, 3).Value = 280
$ws.Cells.Item(90, 4).Value = 419810
$ws.Cells.Item(91, 3).Value = 23
$ws.Cells.Item(91, 4).Value = 34500
$ws.Cells.Item(93, 3).Value = 6278
$ws.Cells.Item(93, 4).Value = 8466938
$ws.Cells.Item(94, 3).Value = 19237
$ws.Cells.Item(94, 4).Value = 27942486
$ws.Cells.Item(95, 3).Value = 44626
$ws.Cells.Item(95, 4).Value = 65843797
$ws.Cells.Item(96, 3).Value = 14274
$ws.Cells.Item(96, 4).Value = 21314949
$ws.Cells.Item(97, 3).Value = 3806
$ws.Cells.Item(97, 4).Value = 5694884
$ws.Cells.Item(98, 3).Value = 652
$ws.Cells.Item(98, 4).Value = 976362
$ws.Cells.Item(99, 3).Value = 31
$ws.Cells.Item(99, 4).Value = 45808
$ws.Cells.Item(101, 3).Value = 16440
$ws.Cells.Item(101, 4).Value = 22347079
$ws.Cells.Item(102, 3).Value = 22004
$ws.Cells.Item(102, 4).Value = 32003066
$ws.Cells.Item(103, 3).Value = 49748
$ws.Cells.Item(103, 4).Value = 73311740
$ws.Cells.Item(104, 3).Value = 15533
$ws.Cells.Item(104, 4).Value = 23167622
$ws.Cells.Item(105, 3).Value = 3979
$ws.Cells.Item(105, 4).Value = 5944814
$ws.Cells.Item(106, 3).Value = 644
$ws.Cells.Item(106, 4).Value = 963054
$ws.Cells.Item(109, 3).Value = 19481
$ws.Cells.Item(109, 4).Value = 26293645
$ws.Cells.Item(110, 3).Value = 8563
$ws.Cells.Item(110, 4).Value = 12510889
$ws.Cells.Item(111, 3).Value = 22178
$ws.Cells.Item(111, 4).Value = 32835213
$ws.Cells.Item(112, 3).Value = 7686
$ws.Cells.Item(112, 4).Value = 11474814
$ws.Cells.Item(113, 3).Value = 1865
$ws.Cells.Item(113, 4).Value = 2790094
$ws.Cells.Item(114, 3).Value = 266
$ws.Cells.Item(114, 4).Value = 397228
$ws.Cells.Item(117, 3).Value = 6970
$ws.Cells.Item(117, 4).Value = 9496898
$ws.Cells.Item(118, 3).Value = 21403
$ws.Cells.Item(118, 4).Value = 31126165
$ws.Cells.Item(119, 3).Value = 52718
$ws.Cells.Item(119, 4).Value = 77770012
$ws.Cells.Item(120, 3).Value = 15865
$ws.Cells.Item(120, 4).Value = 23693664
$ws.Cells.Item(121, 3).Value = 3947
$ws.Cells.Item(121, 4).Value = 5904845
$ws.Cells.Item(122, 3).Value = 804
$ws.Cells.Item(122, 4).Value = 1204212
$ws.Cells.Item(123, 3).Value = 42
$ws.Cells.Item(123, 4).Value = 63000
$ws.Cells.Item(125, 3).Value = 18193
$ws.Cells.Item(125, 4).Value = 25012719
$ws.Cells.Item(126, 3).Value = 29287
$ws.Cells.Item(126, 4).Value = 42904847
$ws.Cells.Item(127, 3).Value = 88098
$ws.Cells.Item(127, 4).Value = 130522924
$ws.Cells.Item(128, 3).Value = 39129
$ws.Cells.Item(128, 4).Value = 58487424
$ws.Cells.Item(129, 3).Value = 12174
$ws.Cells.Item(129, 4).Value = 18234357
$ws.Cells.Item(130, 3).Value = 2467
$ws.Cells.Item(130, 4).Value = 3695665
$ws.Cells.Item(134, 3).Value = 28810
$ws.Cells.Item(134, 4).Value = 40197935
